# "Generate Report for Handoff"
#
# The localization status for e2e\5b2edd21-ede1-464f-9261-522a3c1ef648.md
# moved from "Handed back: in sync with en-US" to "Ready for handoff"
# (because the handback file that came in was not the latest version),
# and the corresponding handoff timestamps / error detail get updated to
# reflect the new handoff that was just generated.

$wb = $excel.ActiveWorkbook

$newStatus        = "Ready for handoff"
$overviewDateTime = "2016-09-06 08:11:19"
$zhHandoffDate    = "2016-09-06 08:10:58"
$deHandoffDate    = "2016-09-06 08:11:19"
$errorDetail      = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80bf7d57d923975860c889700eec8747a68bd0fb/e2e/5b2edd21-ede1-464f-9261-522a3c1ef648.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5c296a71d2f39f40e72c97ade547513e2a943be4/e2e/5b2edd21-ede1-464f-9261-522a3c1ef648.md."

# --- Overview sheet: row for 5b2edd21...md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $overviewDateTime

# --- zh-cn sheet: row for 5b2edd21...md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("H3").Value = $zhHandoffDate
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet: row for 5b2edd21...md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("H3").Value = $deHandoffDate
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
